$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 8344864.5
$ws.Range("I2").Value = 10429330
$ws.Range("J2").Value = 7000
$ws.Range("K2").Value = 10429330
$ws.Range("L2").Value = 7000
$ws.Range("M2").Value = -10429217
$ws.Range("N2").Value = -7226

# Row 18
$ws.Range("H18").Value = 299.08334
$ws.Range("I18").Value = 309
$ws.Range("J18").Value = 190
$ws.Range("K18").Value = 309
$ws.Range("L18").Value = 190
$ws.Range("M18").Value = -25
$ws.Range("N18").Value = -758

# Row 34
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

# Row 36
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

# Row 40
$ws.Range("H40").Value = 2167.6924
$ws.Range("I40").Value = 1756.8572
$ws.Range("J40").Value = 2647
$ws.Range("K40").Value = 1756.8572
$ws.Range("L40").Value = 2647
$ws.Range("M40").Value = -1581.8572
$ws.Range("N40").Value = -2997

# Row 101
$ws.Range("H101").Value = 3097190
$ws.Range("J101").Value = 22727522
$ws.Range("L101").Value = 68182566
$ws.Range("N101").Value = -68185810

# Row 127
$ws.Range("H127").Value = 1108.5714
$ws.Range("I127").Value = 647.2
$ws.Range("J127").Value = 2262
$ws.Range("K127").Value = 1941.6
$ws.Range("L127").Value = 6786
$ws.Range("M127").Value = 3018.4
$ws.Range("N127").Value = -16706


# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2140.1667
$ws.Range("I61").Value = 2230.2856
$ws.Range("J61").Value = 2014
$ws.Range("K61").Value = 2230.2856
$ws.Range("L61").Value = 2014
$ws.Range("M61").Value = -2018.2856
$ws.Range("N61").Value = -2438

# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("N105").Value = 0

# Row 132
$ws.Range("H132").Value = 1728.6333
$ws.Range("I132").Value = 1180.0526
$ws.Range("J132").Value = 2676.182
$ws.Range("K132").Value = 3540.1578
$ws.Range("L132").Value = 8028.545999999999
$ws.Range("M132").Value = -1010.1578
$ws.Range("N132").Value = -13088.546

# Row 136
$ws.Range("H136").Value = 2140.1667
$ws.Range("I136").Value = 2230.2856
$ws.Range("J136").Value = 2014
$ws.Range("K136").Value = 6690.8568
$ws.Range("L136").Value = 6042
$ws.Range("M136").Value = -4140.8568
$ws.Range("N136").Value = -11142

# Row 140
$ws.Range("H140").Value = 41456.855
$ws.Range("J140").Value = 41456.855
$ws.Range("L140").Value = 41456.855
$ws.Range("N140").Value = -51816.855

# Row 141
$ws.Range("H141").Value = 41279.4
$ws.Range("J141").Value = 41279.4
$ws.Range("L141").Value = 41279.4
$ws.Range("N141").Value = -51639.4


# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 657.8946999999999
$ws.Range("I94").Value = 599.2143
$ws.Range("J94").Value = 822.2
$ws.Range("K94").Value = 599.2143
$ws.Range("L94").Value = 822.2
$ws.Range("M94").Value = -148.2143
$ws.Range("N94").Value = -1724.2

# Row 99
$ws.Range("H99").Value = 1755.4762
$ws.Range("I99").Value = 1069.2858
$ws.Range("J99").Value = 2098.5715
$ws.Range("K99").Value = 1069.2858
$ws.Range("L99").Value = 2098.5715
$ws.Range("M99").Value = 428.7141999999999
$ws.Range("N99").Value = -5094.5715

# Row 134
$ws.Range("H134").Value = 4765.706
$ws.Range("I134").Value = 1229.6666
$ws.Range("J134").Value = 8743.75
$ws.Range("K134").Value = 3688.9998
$ws.Range("L134").Value = 26231.25
$ws.Range("M134").Value = -1153.9998
$ws.Range("N134").Value = -31301.25

# Row 140
$ws.Range("H140").Value = 67020
$ws.Range("J140").Value = 67020
$ws.Range("L140").Value = 67020
$ws.Range("N140").Value = -77380


# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1200.9799
$ws.Range("I31").Value = 801.45
$ws.Range("J31").Value = 2883.2104
$ws.Range("K31").Value = 801.45
$ws.Range("L31").Value = 2883.2104
$ws.Range("M31").Value = -506.45
$ws.Range("N31").Value = -3473.2104

# Row 34
$ws.Range("H34").Value = 1200.9799
$ws.Range("I34").Value = 801.45
$ws.Range("J34").Value = 2883.2104
$ws.Range("K34").Value = 801.45
$ws.Range("L34").Value = 2883.2104
$ws.Range("M34").Value = -599.45
$ws.Range("N34").Value = -3287.2104

# Row 51
$ws.Range("H51").Value = 9800
$ws.Range("I51").Value = 9800
$ws.Range("K51").Value = 9800
$ws.Range("M51").Value = -9064

# Row 59
$ws.Range("H59").Value = 27100
$ws.Range("I59").Value = 5000
$ws.Range("J59").Value = 31520
$ws.Range("K59").Value = 5000
$ws.Range("L59").Value = 31520
$ws.Range("M59").Value = -3855
$ws.Range("N59").Value = -33810

# Row 60
$ws.Range("H60").Value = 24304.87
$ws.Range("I60").Value = 2250
$ws.Range("J60").Value = 26405.334
$ws.Range("K60").Value = 2250
$ws.Range("L60").Value = 26405.334
$ws.Range("M60").Value = -1739
$ws.Range("N60").Value = -27427.334

# Row 61
$ws.Range("H61").Value = 9800
$ws.Range("I61").Value = 9800
$ws.Range("K61").Value = 9800
$ws.Range("M61").Value = -9452

# Row 68
$ws.Range("H68").Value = 22795
$ws.Range("J68").Value = 22795
$ws.Range("L68").Value = 22795
$ws.Range("N68").Value = -24293

# Row 71
$ws.Range("H71").Value = 22795
$ws.Range("J71").Value = 22795
$ws.Range("L71").Value = 68385
$ws.Range("N71").Value = -75873

# Row 74
$ws.Range("H74").Value = 5285
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("H77").Value = 5285
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# Row 134
$ws.Range("H134").Value = 860.1923
$ws.Range("I134").Value = 834.6
$ws.Range("K134").Value = 2503.8
$ws.Range("M134").Value = 31.19999999999982


# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 140
$ws.Range("I86").Value = 100
$ws.Range("K86").Value = 300
$ws.Range("M86").Value = 886

# Row 89
$ws.Range("H89").Value = 140
$ws.Range("I89").Value = 100
$ws.Range("K89").Value = 900
$ws.Range("M89").Value = 5028


# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 2904.25
$ws.Range("I43").Value = 2904.25
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 2904.25
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -2753.25

# Row 57
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").ClearContents()
$ws.Range("N57").Value = 0

# Row 80
$ws.Range("H80").Value = 12504674
$ws.Range("I80").Value = 8466.333000000001
$ws.Range("J80").Value = 20002398
$ws.Range("K80").Value = 8466.333000000001
$ws.Range("L80").Value = 20002398
$ws.Range("M80").Value = -7468.333000000001
$ws.Range("N80").Value = -20004394

# Row 83
$ws.Range("H83").Value = 12504674
$ws.Range("I83").Value = 8466.333000000001
$ws.Range("J83").Value = 20002398
$ws.Range("K83").Value = 42331.665
$ws.Range("L83").Value = 100011990
$ws.Range("M83").Value = -37339.665
$ws.Range("N83").Value = -100021974

# Row 132
$ws.Range("H132").Value = 4991.1
$ws.Range("I132").Value = 1400.3334
$ws.Range("J132").Value = 13369.556
$ws.Range("K132").Value = 4201.0002
$ws.Range("L132").Value = 40108.66800000001
$ws.Range("M132").Value = -1671.0002
$ws.Range("N132").Value = -45168.66800000001

# Row 138
$ws.Range("H138").Value = 60381.08
$ws.Range("J138").Value = 60381.08
$ws.Range("L138").Value = 60381.08
$ws.Range("N138").Value = -70661.08


# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1363.7407
$ws.Range("I7").Value = 1262.4286
$ws.Range("J7").Value = 1718.3334
$ws.Range("K7").Value = 1262.4286
$ws.Range("L7").Value = 1718.3334
$ws.Range("M7").Value = -1150.4286
$ws.Range("N7").Value = -1942.3334

# Row 122
$ws.Range("H122").Value = 5708.6763
$ws.Range("I122").Value = 6369.423
$ws.Range("K122").Value = 19108.269
$ws.Range("M122").Value = -16658.269

# Row 126
$ws.Range("H126").Value = 1363.7407
$ws.Range("I126").Value = 1262.4286
$ws.Range("J126").Value = 1718.3334
$ws.Range("K126").Value = 3787.2858
$ws.Range("L126").Value = 5155.0002
$ws.Range("M126").Value = -1317.2858
$ws.Range("N126").Value = -10095.0002

# Row 138
$ws.Range("H138").Value = 33441.8
$ws.Range("J138").Value = 33441.8
$ws.Range("L138").Value = 33441.8
$ws.Range("N138").Value = -43721.8

# Row 140
$ws.Range("H140").Value = 55636.375
$ws.Range("J140").Value = 55636.375
$ws.Range("L140").Value = 55636.375
$ws.Range("N140").Value = -65996.375


# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1239.0385
$ws.Range("I96").Value = 947.94116
$ws.Range("J96").Value = 1788.8889
$ws.Range("K96").Value = 947.94116
$ws.Range("L96").Value = 1788.8889
$ws.Range("M96").Value = 425.05884
$ws.Range("N96").Value = -4534.8889

# Row 126
$ws.Range("H126").Value = 582.9231
$ws.Range("I126").Value = 470.03125
$ws.Range("K126").Value = 1410.09375
$ws.Range("M126").Value = 1059.90625

# Row 135
$ws.Range("H135").Value = 45000
$ws.Range("J135").Value = 45000
$ws.Range("L135").Value = 45000
$ws.Range("N135").Value = -55140

# Row 137
$ws.Range("H137").Value = 69000
$ws.Range("J137").Value = 69000
$ws.Range("L137").Value = 69000
$ws.Range("N137").Value = -79200

# Row 140
$ws.Range("H140").Value = 70701.5
$ws.Range("J140").Value = 70701.5
$ws.Range("L140").Value = 70701.5
$ws.Range("N140").Value = -81061.5

